$d = $word.ActiveDocument

# Locate the paragraph that introduces the Perseus constellation
# ("... des Sternbildes Perseus am Nachthimmel ...").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Sternbildes*Perseus*Nachthimmel*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the Perseus constellation paragraph"
}

$newText = "Mach mit an einer weltweiten Kampagne, die schwächsten sichtbaren Sterne zu beobachten und aufzuzeichnen, um die Lichtverschmutzung an einem Ort zu messen. Durch das Auffinden und Beobachten des Perseus-Konstellation am Nachthimmel und den Vergleich mit den Helligkeitskarten, lernen Menschen auf der ganzen Erde, wie die Lichter in ihrer Gemeinde zur Lichtverschmutzung beitragen. Dein Beitrag zur Online-Datenbank beschreibt den sichtbaren Nachthimmel."

# Replace the whole run content (but keep the paragraph mark / pPr intact)
# with a single freshly-inserted run that carries no explicit run
# formatting, collapsing the four previously separate runs into one.
$body = $d.Range($target.Range.Start, $target.Range.End - 1)
$body.Delete()
$target.Range.InsertBefore($newText)
